# ------------------------------------------------------------------
# Rebuild Language.xlsx: Sheet1 -> "Comm" plus four new localization
# sheets (Property, Guild, Tip, Item), matching the target layout.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ====================================================================
# 1) "Comm" sheet (the original, renamed)
# ====================================================================
$comm = $wb.Worksheets.Item(1)
$comm.Name = "Comm"

# Column A: localization keys (was 中文_N -> now Langage_Comm_N)
$comm.Range("A2").Value = "Langage_Comm_1"
$comm.Range("A3").Value = "Langage_Comm_2"
$comm.Range("A4").Value = "Langage_Comm_3"
$comm.Range("A5").Value = "Langage_Comm_4"
$comm.Range("A6").Value = "Langage_Comm_5"
$comm.Range("A7").Value = "Langage_Comm_6"

# Column C: the Chinese text (new wording)
$comm.Range("C2").Value = "确认"
$comm.Range("C3").Value = "取消"
$comm.Range("C4").Value = "登录"
$comm.Range("C5").Value = "创建角色"
$comm.Range("C6").Value = "进入游戏"
$comm.Range("C7").Value = "中文_6"

# New blank rows 8-12, formatted like row 7 (style s="2")
$comm.Range("A7:C7").Copy()
$comm.Range("A8:C12").PasteSpecial(-4122)

# Column widths (engine snaps to pixel grid - closest achievable to
# the 31.875 / 24.5 / 23 target)
$comm.Columns.Item(1).ColumnWidth = 31.142857142857142
$comm.Columns.Item(2).ColumnWidth = 23.714285714285715
$comm.Columns.Item(3).ColumnWidth = 22.285714285714285

$comm.Range("C8").Select()

# ====================================================================
# 2) "Property" sheet - fresh sheet listing character stats
#    (added now, right after Comm; "Guild" is spliced in between the
#    two further down so the final tab order is Comm / Property /
#    Guild / Tip / Item)
# ====================================================================
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$property = $wb.Worksheets.Add($null, $last)
$property.Name = "Property"

$comm.Range("A1:C1").Copy()
$property.Range("A1:C1").PasteSpecial(-4163)
$comm.Range("A1:C1").Copy()
$property.Range("A1:C1").PasteSpecial(-4122)

$property.Range("A2").Value = "Langage_HP"
$property.Range("A3").Value = "Langage_MAXHP"
$property.Range("A4").Value = "Langage_MP"
$property.Range("A5").Value = "Langage_MAXMP"
$property.Range("A6").Value = "Langage_VP"
$property.Range("A7").Value = "Langage_ATTACK"

$property.Range("A2:A7").Copy()
$property.Range("A8:A28").PasteSpecial(-4122)
$property.Range("A8:A28").ClearContents()

$property.Columns.Item(1).ColumnWidth = 45.142857142857146

$property.Range("A1:XFD1").Select()

# ====================================================================
# 4) "Tip" sheet - fresh sheet, header row only
# ====================================================================
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$tip = $wb.Worksheets.Add($null, $last)
$tip.Name = "Tip"

$comm.Range("A1:C1").Copy()
$tip.Range("A1:C1").PasteSpecial(-4163)
$comm.Range("A1:C1").Copy()
$tip.Range("A1:C1").PasteSpecial(-4122)

$tip.Range("A1:XFD1").Select()

# ====================================================================
# 5) "Item" sheet - fresh sheet, header row only
# ====================================================================
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$item = $wb.Worksheets.Add($null, $last)
$item.Name = "Item"

$comm.Range("A1:C1").Copy()
$item.Range("A1:C1").PasteSpecial(-4163)
$comm.Range("A1:C1").Copy()
$item.Range("A1:C1").PasteSpecial(-4122)

$item.Range("A1:XFD1").Select()

# ====================================================================
# Finish on the Comm tab
# ====================================================================
$comm.Select()
